$wb = $excel.ActiveWorkbook
$aw = $excel.ActiveWindow
try { $aw.WindowWidth = 25280; Write-Output "set WindowWidth ok" } catch { Write-Output "WindowWidth failed: $_" }
try { $aw.WindowHeight = 26420; Write-Output "set WindowHeight ok" } catch { Write-Output "WindowHeight failed: $_" }
